# Actualización automática 2025-11-28 16:30:09
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Widen column J (10th column) from 10 to 11 characters
$ws1.Range("J1").ColumnWidth = 10.15

$ws1.Range("M3").Value = 498.96

$ws1.Range("J4").Value = 12.29
$ws1.Range("M4").Value = 9466.24

$ws1.Range("H7").Value = 837

$ws1.Range("L16").Value = 4752
$ws1.Range("M16").Value = 14106.02

$ws1.Range("M18").Value = 2383.22

$ws1.Range("M41").Value = 641.5

$ws1.Range("M53").Value = 328.86

$ws1.Range("M56").Value = "19 de 54"

# ---------------------------------------------------------------------------
# Sheet 2: VENTA MENSUAL
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("F3").Value = 498.96
$ws2.Range("F4").Value = 12438.61
$ws2.Range("F7").Value = 702.21
$ws2.Range("F16").Value = 19528.93
$ws2.Range("F18").Value = 3299.06
$ws2.Range("F41").Value = 641.5
$ws2.Range("F55").Value = 1011.96
$ws2.Range("F56").Value = 1011.96
$ws2.Range("F60").Value = 102709.52

# ---------------------------------------------------------------------------
# Sheet 3: CUMPLIMIENTO MENSUAL
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Range("D6").Value = 4611.11
$ws3.Range("E6").Value = -1703.52631853974
$ws3.Range("F6").Value = 1.58589072754879

$ws3.Range("D9").Value = 165.79
$ws3.Range("E9").Value = -165.79

$ws3.Range("D11").Value = 23642.66
$ws3.Range("E11").Value = -9406.67
$ws3.Range("F11").Value = 1.660766831109041

$ws3.Range("D12").Value = 55522.16
$ws3.Range("E12").Value = 9421.839999999997
$ws3.Range("F12").Value = 0.8549236265089925

$ws3.Range("D14").Value = 99218.35000000001
$ws3.Range("E14").Value = -262.0931407616256
$ws3.Range("F14").Value = 1.00264857573518
